$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new value, whether the cell holds a numeric-looking
# string that must stay text (column D prices use "." as a thousands separator,
# so Excel would otherwise silently reinterpret them as numbers).
$updates = @(
    ,@("D2", "63.486.87")
    ,@("E2", "  +0.62%  ")
    ,@("D3", "3.097.76")
    ,@("E3", "  -0.50%  ")
    ,@("E4", "  -0.03%  ")
    ,@("D5", "584.15")
    ,@("E5", "  -0.14%  ")
    ,@("D6", "144.80")
    ,@("E6", "  +0.27%  ")
    ,@("E7", "  +0.01%  ")
    ,@("D8", "3.091.89")
    ,@("E8", "  -0.45%  ")
    ,@("E9", "  -0.27%  ")
    ,@("D10", "0.160")
    ,@("E10", "  +6.69%  ")
    ,@("E11", "  -2.93%  ")
    ,@("D12", "0.457")
    ,@("E12", "  -2.07%  ")
    ,@("D13", "0.0000245")
    ,@("E13", "  -0.96%  ")
    ,@("D14", "37.02")
    ,@("E14", "  +4.26%  ")
    ,@("E15", "  -1.11%  ")
    ,@("D16", "3.611.52")
    ,@("E16", "  -0.51%  ")
    ,@("D17", "63.361.13")
    ,@("E17", "  +0.56%  ")
    ,@("B18", "WrappedEther")
    ,@("C18", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth")
    ,@("D18", "3.097.48")
    ,@("E18", "  -0.50%  ")
    ,@("B19", "Polkadot")
    ,@("C19", "https://coinranking.com/coin/25W7FG7om+polkadot-dot")
    ,@("D19", "7.06")
    ,@("E19", "  -1.37%  ")
    ,@("D20", "459.83")
    ,@("E20", "  -1.58%  ")
    ,@("D21", "14.23")
    ,@("E21", "  +1.01%  ")
    ,@("E22", "  -0.48%  ")
    ,@("E23", "  -1.73%  ")
    ,@("D24", "12.92")
    ,@("E24", "  -2.64%  ")
    ,@("D25", "81.02")
    ,@("D26", "2.24")
    ,@("E26", "  +3.23%  ")
    ,@("E27", "  +0.00%  ")
    ,@("D28", "9.22")
    ,@("E28", "  +9.98%  ")
    ,@("B29", "FirstDigitalUSD")
    ,@("C29", "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd")
    ,@("D29", "1.00")
    ,@("E29", "  +0.00%  ")
    ,@("B30", "PancakeSwap")
    ,@("C30", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake")
    ,@("D30", "2.67")
    ,@("E30", "  -0.14%  ")
    ,@("E31", "  -1.09%  ")
    ,@("D32", "6.95")
    ,@("E32", "  +1.79%  ")
    ,@("E33", "  +1.00%  ")
    ,@("D34", "26.62")
    ,@("E34", "  -1.02%  ")
    ,@("D35", "0.0₃0844")
    ,@("E35", "  -1.91%  ")
    ,@("E36", "  -0.53%  ")
    ,@("B37", "Stacks")
    ,@("C37", "https://coinranking.com/coin/mMPrMcB7+stacks-stx")
    ,@("D37", "2.30")
    ,@("E37", "  -4.51%  ")
    ,@("B38", "dogwifhat")
    ,@("C38", "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif")
    ,@("D38", "3.35")
    ,@("E38", "  +1.44%  ")
    ,@("E39", "  -0.49%  ")
    ,@("E40", "  -1.26%  ")
    ,@("D41", "434.97")
    ,@("E41", "  +1.13%  ")
    ,@("E42", "  -0.35%  ")
    ,@("E43", "  -0.24%  ")
    ,@("D44", "2.882.20")
    ,@("E44", "  -1.47%  ")
    ,@("D45", "0.275")
    ,@("E45", "  -1.43%  ")
    ,@("E46", "  -2.75%  ")
    ,@("D47", "36.33")
    ,@("E47", "  +2.83%  ")
    ,@("D48", "125.56")
    ,@("E48", "  +1.83%  ")
    ,@("E49", "  -0.01%  ")
    ,@("E50", "  -1.13%  ")
    ,@("D51", "24.09")
    ,@("E51", "  -1.66%  ")
)

foreach ($u in $updates) {
    $addr = $u[0]
    $val = $u[1]
    $col = $addr.Substring(0, 1)
    $rng = $ws.Range($addr)
    if ($col -eq "D") {
        # Force text so numeric-looking price strings (e.g. "584.15",
        # "63.486.87") are preserved verbatim instead of becoming numbers.
        $rng.NumberFormat = "@"
        $rng.Value = $val
        $rng.Style = "Normal"
    } else {
        $rng.Value = $val
    }
}
